$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.0507701114283492
$ws.Range("D2").Value = 0.0002214791267642902
$ws.Range("E2").Value = 0.4330778936147368
$ws.Range("F2").Value = 0.6071500647315276
$ws.Range("G2").Value = 0.5211106091079785
$ws.Range("H2").Value = 0.4720876509364018
$ws.Range("N2").Value = 2.873581687696344
$ws.Range("O2").Value = 1.947170118530437
# Row 3
$ws.Range("C3").Value = 0.0449891660834254
$ws.Range("D3").Value = 0.0002074950308945667
$ws.Range("E3").Value = 0.3774308521777527
$ws.Range("F3").Value = 0.5604074233862946
$ws.Range("G3").Value = 0.4694198853888736
$ws.Range("H3").Value = 0.4524053179119107
$ws.Range("N3").Value = 2.562605684679454
$ws.Range("O3").Value = 1.795960797759165
# Row 4
$ws.Range("C4").Value = 0.04145654160221568
$ws.Range("D4").Value = 0.0001990096585942069
$ws.Range("E4").Value = 0.3433846769561768
$ws.Range("F4").Value = 0.5320651601840893
$ws.Range("G4").Value = 0.4379383000529629
$ws.Range("H4").Value = 0.4406200446956348
$ws.Range("N4").Value = 2.371325805375761
$ws.Range("O4").Value = 1.704272837477333
# Row 5
$ws.Range("C5").Value = 0.04002112976345984
$ws.Range("D5").Value = 0.0001955709880077094
$ws.Range("E5").Value = 0.3295384505563845
$ws.Range("F5").Value = 0.5206045375284276
$ws.Range("G5").Value = 0.4251727262033285
$ws.Range("H5").Value = 0.4358924360248295
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 1.667196503579476
# Row 6
$ws.Range("C6").Value = 0.03978302949464307
$ws.Range("D6").Value = 0.000195000981014104
$ws.Range("E6").Value = 0.327240901457813
$ws.Range("F6").Value = 0.5187068634029828
$ws.Range("G6").Value = 0.4230568104949839
$ws.Range("H6").Value = 0.4351119378590909
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 1.661057262505437
# Row 7
$ws.Range("C7").Value = 0.04143716640180628
$ws.Range("D7").Value = 0.0001989632145475984
$ws.Range("E7").Value = 0.3431978325289862
$ws.Range("F7").Value = 0.5319102387528574
$ws.Range("G7").Value = 0.4377658837299805
$ws.Range("H7").Value = 0.4405559834186477
$ws.Range("N7").Value = 2.370273851395496
$ws.Range("O7").Value = 1.703771654249749
# Row 8
$ws.Range("C8").Value = 0.04877327908289431
$ws.Range("D8").Value = 0.0002166326728543222
$ws.Range("E8").Value = 0.4138639980941292
$ws.Range("F8").Value = 0.5909583477640865
$ws.Range("G8").Value = 0.5032336662155217
$ws.Range("H8").Value = 0.4652387236575066
$ws.Range("N8").Value = 2.766433886209654
$ws.Range("O8").Value = 1.894791391600563
# Row 9
$ws.Range("C9").Value = 0.06329817106426106
$ws.Range("D9").Value = 0.0002523725686081058
$ws.Range("E9").Value = 0.5535289322776578
$ws.Range("F9").Value = 0.7096393624568123
$ws.Range("G9").Value = 0.6337105273447037
$ws.Range("H9").Value = 0.5160408099667677
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 2.278712888236214
# Row 10
$ws.Range("C10").Value = 0.07406179888684505
$ws.Range("D10").Value = 0.000279711252257675
$ws.Range("E10").Value = 0.6569889858713083
$ws.Range("F10").Value = 0.7986707726857531
$ws.Range("G10").Value = 0.7309386318265467
$ws.Range("H10").Value = 0.5548613196053509
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 2.566730516189295
# Row 11
$ws.Range("C11").Value = 0.07898017322646922
$ws.Range("D11").Value = 0.0002924779473816308
$ws.Range("E11").Value = 0.7042806635774355
$ws.Range("F11").Value = 0.8395893270043473
$ws.Range("G11").Value = 0.7754870784960133
$ws.Range("H11").Value = 0.5728542709317139
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 2.699107961274422
# Row 12
$ws.Range("C12").Value = 0.08084590513031742
$ws.Range("D12").Value = 0.0002973680175808457
$ws.Range("E12").Value = 0.7222246361523901
$ws.Range("F12").Value = 0.8551453888328808
$ws.Range("G12").Value = 0.7924037097205883
$ws.Range("H12").Value = 0.5797161956365926
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 2.749435134041391
# Row 13
$ws.Range("C13").Value = 0.08044394031165325
$ws.Range("D13").Value = 0.0002963122483514713
$ws.Range("E13").Value = 0.71835844768475
$ws.Range("F13").Value = 0.8517923789561905
$ws.Range("G13").Value = 0.788758295682527
$ws.Range("H13").Value = 0.5782361981204076
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 2.738587381443438
# Row 14
$ws.Range("C14").Value = 0.0791336024293372
$ws.Range("D14").Value = 0.0002928790944050075
$ws.Range("E14").Value = 0.705756193525616
$ws.Range("F14").Value = 0.8408679032461208
$ws.Range("G14").Value = 0.776877868763421
$ws.Range("H14").Value = 0.5734178332718614
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 2.703244404860413
# Row 15
$ws.Range("C15").Value = 0.07833140874255662
$ws.Range("D15").Value = 0.0002907836756118609
$ws.Range("E15").Value = 0.6980416793989974
$ws.Range("F15").Value = 0.8341843332764398
$ws.Range("G15").Value = 0.76960693541767
$ws.Range("H15").Value = 0.5704727604213815
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 2.681621796540469
# Row 16
$ws.Range("C16").Value = 0.07374082122164793
$ws.Range("D16").Value = 0.0002788842498977484
$ws.Range("E16").Value = 0.653903188055736
$ws.Range("F16").Value = 0.7960051418863969
$ws.Range("G16").Value = 0.7280338090902205
$ws.Range("H16").Value = 0.5536921836557838
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 2.558106947690987
# Row 17
$ws.Range("C17").Value = 0.07093033230552237
$ws.Range("D17").Value = 0.0002716745862141323
$ws.Range("E17").Value = 0.6268859407979761
$ws.Range("F17").Value = 0.7726911580841289
$ws.Range("G17").Value = 0.7026125999634871
$ws.Range("H17").Value = 0.5434835327375822
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 2.482684641670062
# Row 18
$ws.Range("C18").Value = 0.06931587337896872
$ws.Range("D18").Value = 0.0002675583527864678
$ws.Range("E18").Value = 0.6113675376912369
$ws.Range("F18").Value = 0.7593208267845597
$ws.Range("G18").Value = 0.6880209333433243
$ws.Range("H18").Value = 0.5376431314455488
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 2.439431214395086
# Row 19
$ws.Range("C19").Value = 0.0687695962650281
$ws.Range("D19").Value = 0.0002661696763892607
$ws.Range("E19").Value = 0.6061168248240705
$ws.Range("F19").Value = 0.7548005734839478
$ws.Range("G19").Value = 0.6830855503779674
$ws.Range("H19").Value = 0.5356710416518524
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 2.424808134334342
# Row 20
$ws.Range("C20").Value = 0.07122929976881665
$ws.Range("D20").Value = 0.0002724388452903792
$ws.Range("E20").Value = 0.6297597570736997
$ws.Range("F20").Value = 0.7751689012727923
$ws.Range("G20").Value = 0.7053156233760376
$ws.Range("H20").Value = 0.5445670137242189
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("O20").Value = 2.490700254398519
# Row 21
$ws.Range("C21").Value = 0.07951839145849249
$ws.Range("D21").Value = 0.0002938859195724675
$ws.Range("E21").Value = 0.7094567898753752
$ws.Range("F21").Value = 0.8440750202745733
$ws.Range("G21").Value = 0.7803661498241468
$ws.Range("H21").Value = 0.574831787014972
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("O21").Value = 2.713620071985645
# Row 22
$ws.Range("C22").Value = 0.08495478170114268
$ws.Range("D22").Value = 0.0003082308938353595
$ws.Range("E22").Value = 0.7617525596761112
$ws.Range("F22").Value = 0.8894657372637198
$ws.Range("G22").Value = 0.8296911384626071
$ws.Range("H22").Value = 0.5948938050911465
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("O22").Value = 2.860470985179404
# Row 23
$ws.Range("C23").Value = 0.08205150890684365
$ws.Range("D23").Value = 0.0003005419286825628
$ws.Range("E23").Value = 0.7338211921436084
$ws.Range("F23").Value = 0.8652068872279415
$ws.Range("G23").Value = 0.8033398739482038
$ws.Range("H23").Value = 0.5841603510018558
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 2.781986552202113
# Row 24
$ws.Range("C24").Value = 0.07109413244812401
$ws.Range("D24").Value = 0.0002720932351916172
$ws.Range("E24").Value = 0.6284604608163988
$ws.Range("F24").Value = 0.7740486102795359
$ws.Range("G24").Value = 0.7040935139324915
$ws.Range("H24").Value = 0.5440770826470498
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 2.487076060191782
# Row 25
$ws.Range("C25").Value = 0.0593531637428697
$ws.Range("D25").Value = 0.0002425453255574528
$ws.Range("E25").Value = 0.5156103763282829
$ws.Range("F25").Value = 0.6772160331683352
$ws.Range("G25").Value = 0.5981795963497234
$ws.Range("H25").Value = 0.5020371973616591
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 2.173826014020051
